$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = -7.670899999999997
$ws.Range("B3").Value = 5.846100000000003
$ws.Range("E3").Value = 15.63659999999999
$ws.Range("E12").Value = 17.44740000000001
$ws.Range("B14").Value = 5.586099999999999
$ws.Range("B16").Value = 6.174400000000003
$ws.Range("D18").Value = -8.935699999999999
$ws.Range("B21").Value = 9.129700000000003
$ws.Range("B23").Value = 9.2623
$ws.Range("D24").Value = -7.305200000000005
$ws.Range("E24").Value = 16.59930000000001
$ws.Range("B25").Value = 5.731399999999997
$ws.Range("D25").Value = -7.465099999999995
$ws.Range("E25").Value = 17.0979
$ws.Range("B26").Value = 6.701600000000002
$ws.Range("D27").Value = -8.900700000000008
$ws.Range("B29").Value = 4.791499999999997
$ws.Range("D30").Value = -7.708600000000005
$ws.Range("D31").Value = -8.533300000000008
$ws.Range("D39").Value = -8.187399999999998
$ws.Range("B40").Value = 8.7623
$ws.Range("E41").Value = 16.2326
$ws.Range("D42").Value = -7.983000000000001
$ws.Range("D48").Value = -7.326499999999998
$ws.Range("E50").Value = 16.43780000000001
$ws.Range("D51").Value = -7.755199999999999
$ws.Range("D52").Value = -7.4631
$ws.Range("B53").Value = 5.398000000000001
$ws.Range("E53").Value = 16.517
$ws.Range("D55").Value = -8.843200000000001
$ws.Range("D56").Value = -7.974099999999998
$ws.Range("E56").Value = 16.72300000000001
$ws.Range("B57").Value = 4.755499999999999
$ws.Range("D57").Value = -8.222699999999998
$ws.Range("E57").Value = 16.48200000000001
$ws.Range("E58").Value = 16.13830000000002
$ws.Range("B59").Value = 6.105599999999997
$ws.Range("D60").Value = -7.756799999999997
$ws.Range("E61").Value = 16.48730000000001
$ws.Range("E63").Value = 17.49030000000001
$ws.Range("E64").Value = 17.4688
$ws.Range("B65").Value = 6.122200000000001
$ws.Range("B69").Value = 5.138899999999995
$ws.Range("E70").Value = 17.34400000000001
$ws.Range("E72").Value = 16.972
$ws.Range("D73").Value = -8.136599999999994
$ws.Range("D74").Value = -7.679100000000009
$ws.Range("B79").Value = 9.463300000000006
$ws.Range("B83").Value = 4.620799999999999
$ws.Range("E86").Value = 16.64260000000001
$ws.Range("D89").Value = -7.137299999999994
$ws.Range("E89").Value = 17.43440000000002
$ws.Range("D90").Value = -8.139800000000003
$ws.Range("B91").Value = 5.4095
$ws.Range("D92").Value = -5.773100000000003
$ws.Range("B93").Value = 5.955099999999998
$ws.Range("E98").Value = 15.21810000000001
$ws.Range("B100").Value = 5.579700000000002
$ws.Range("E100").Value = 16.89340000000001
$ws.Range("E102").Value = 16.81179999999998
